$wb = $excel.ActiveWorkbook
$users = $wb.Worksheets.Item(1)
$incomes = $wb.Worksheets.Item(2)

# --- users sheet: add a new "role" column definition row (row 8) ---

# Clone the formatting (borders/fonts/fills) of row 7 into the new row 8
# so the new cells reuse the existing style entries instead of minting new
# ones (xlPasteFormats = -4122).
$users.Range("A7:G7").Copy()
$users.Range("A8:G8").PasteSpecial(-4122)

$users.Range("A8").Value = 7
$users.Range("B8").Value = "role"
$users.Range("C8").Value = "権限"
$users.Range("D8").Value = "VARCHAR（10）"
$users.Range("F8").Value = "×"

# Existing row 5 "Null" marker re-saved (picks up the phonetic-guide variant
# of the "×" string already used elsewhere in the sheet).
$users.Range("F5").Value = "×"

# --- view state: "users" tab becomes the active/selected sheet ---
$users.Activate() | Out-Null
$users.Range("F8").Select() | Out-Null
